$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7676553333333332
$ws.Cells.Item(2, 8).Value = 2.302966
$ws.Cells.Item(2, 9).Value = 0.3736977786965754
$ws.Cells.Item(2, 10).Value = 0.3736977786965754
$ws.Cells.Item(2, 13).Value = 2.318119
$ws.Cells.Item(2, 14).Value = 6.954357
$ws.Cells.Item(2, 15).Value = 0.1070970465647729
$ws.Cells.Item(2, 16).Value = 0.1070970465647729
$ws.Cells.Item(2, 17).Value = 1.779516413651333
$ws.Cells.Item(2, 18).Value = 16.015647722862
$ws.Cells.Item(2, 19).Value = 0.04002192840621934
$ws.Cells.Item(2, 20).Value = 0.04002192840621934

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7676553333333332
$ws.Cells.Item(3, 8).Value = 2.302966
$ws.Cells.Item(3, 9).Value = 0.3736977786965754
$ws.Cells.Item(3, 10).Value = 0.3736977786965754
$ws.Cells.Item(3, 14).Value = 5.768654000000001
$ws.Cells.Item(3, 15).Value = 0.0888372291002696
$ws.Cells.Item(3, 16).Value = 0.0888372291002696
$ws.Cells.Item(3, 17).Value = 1.476112669751555
$ws.Cells.Item(3, 18).Value = 13.285014027764
$ws.Cells.Item(3, 19).Value = 0.03319827518032952
$ws.Cells.Item(3, 20).Value = 0.03319827518032951

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.7676553333333332
$ws.Cells.Item(4, 8).Value = 2.302966
$ws.Cells.Item(4, 9).Value = 0.3736977786965754
$ws.Cells.Item(4, 10).Value = 0.3736977786965754
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.108471
$ws.Cells.Item(4, 14).Value = 0.325413
$ws.Cells.Item(4, 15).Value = 0.005011357802566427
$ws.Cells.Item(4, 16).Value = 0.005011357802566428
$ws.Cells.Item(4, 17).Value = 0.08326834166199999
$ws.Cells.Item(4, 18).Value = 0.7494150749579999
$ws.Cells.Item(4, 19).Value = 0.001872733279072825
$ws.Cells.Item(4, 20).Value = 0.001872733279072825

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.7676553333333332
$ws.Cells.Item(5, 8).Value = 2.302966
$ws.Cells.Item(5, 9).Value = 0.3736977786965754
$ws.Cells.Item(5, 10).Value = 0.3736977786965754
$ws.Cells.Item(5, 13).Value = 17.29555733333333
$ws.Cells.Item(5, 14).Value = 51.886672
$ws.Cells.Item(5, 15).Value = 0.7990543665323911
$ws.Cells.Item(5, 16).Value = 0.7990543665323911
$ws.Cells.Item(5, 17).Value = 13.27702682990578
$ws.Cells.Item(5, 18).Value = 119.493241469152
$ws.Cells.Item(5, 19).Value = 0.2986048418309538
$ws.Cells.Item(5, 20).Value = 0.2986048418309538

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.811404
$ws.Cells.Item(6, 8).Value = 2.434212
$ws.Cells.Item(6, 9).Value = 0.3949948098567449
$ws.Cells.Item(6, 10).Value = 0.3949948098567449
$ws.Cells.Item(6, 13).Value = 2.318119
$ws.Cells.Item(6, 14).Value = 6.954357
$ws.Cells.Item(6, 15).Value = 0.1070970465647729
$ws.Cells.Item(6, 16).Value = 0.1070970465647729
$ws.Cells.Item(6, 17).Value = 1.880931029076
$ws.Cells.Item(6, 18).Value = 16.928379261684
$ws.Cells.Item(6, 19).Value = 0.04230277754407143
$ws.Cells.Item(6, 20).Value = 0.04230277754407143

$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.811404
$ws.Cells.Item(7, 8).Value = 2.434212
$ws.Cells.Item(7, 9).Value = 0.3949948098567449
$ws.Cells.Item(7, 10).Value = 0.3949948098567449
$ws.Cells.Item(7, 14).Value = 5.768654000000001
$ws.Cells.Item(7, 15).Value = 0.0888372291002696
$ws.Cells.Item(7, 16).Value = 0.0888372291002696
$ws.Cells.Item(7, 17).Value = 1.560236310072
$ws.Cells.Item(7, 18).Value = 14.042126790648
$ws.Cells.Item(7, 19).Value = 0.03509024441666107
$ws.Cells.Item(7, 20).Value = 0.03509024441666107

$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.811404
$ws.Cells.Item(8, 8).Value = 2.434212
$ws.Cells.Item(8, 9).Value = 0.3949948098567449
$ws.Cells.Item(8, 10).Value = 0.3949948098567449
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.108471
$ws.Cells.Item(8, 14).Value = 0.325413
$ws.Cells.Item(8, 15).Value = 0.005011357802566427
$ws.Cells.Item(8, 16).Value = 0.005011357802566428
$ws.Cells.Item(8, 17).Value = 0.088013803284
$ws.Cells.Item(8, 18).Value = 0.792124229556
$ws.Cells.Item(8, 19).Value = 0.00197946032234884
$ws.Cells.Item(8, 20).Value = 0.001979460322348841

$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.811404
$ws.Cells.Item(9, 8).Value = 2.434212
$ws.Cells.Item(9, 9).Value = 0.3949948098567449
$ws.Cells.Item(9, 10).Value = 0.3949948098567449
$ws.Cells.Item(9, 13).Value = 17.29555733333333
$ws.Cells.Item(9, 14).Value = 51.886672
$ws.Cells.Item(9, 15).Value = 0.7990543665323911
$ws.Cells.Item(9, 16).Value = 0.7990543665323911
$ws.Cells.Item(9, 17).Value = 14.033684402496
$ws.Cells.Item(9, 18).Value = 126.303159622464
$ws.Cells.Item(9, 19).Value = 0.3156223275736636
$ws.Cells.Item(9, 20).Value = 0.3156223275736636

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.475155
$ws.Cells.Item(10, 8).Value = 1.425465
$ws.Cells.Item(10, 9).Value = 0.2313074114466796
$ws.Cells.Item(10, 10).Value = 0.2313074114466796
$ws.Cells.Item(10, 13).Value = 2.318119
$ws.Cells.Item(10, 14).Value = 6.954357
$ws.Cells.Item(10, 15).Value = 0.1070970465647729
$ws.Cells.Item(10, 16).Value = 0.1070970465647729
$ws.Cells.Item(10, 17).Value = 1.101465833445
$ws.Cells.Item(10, 18).Value = 9.913192501005
$ws.Cells.Item(10, 19).Value = 0.02477234061448213
$ws.Cells.Item(10, 20).Value = 0.02477234061448214

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.475155
$ws.Cells.Item(11, 8).Value = 1.425465
$ws.Cells.Item(11, 9).Value = 0.2313074114466796
$ws.Cells.Item(11, 10).Value = 0.2313074114466796
$ws.Cells.Item(11, 14).Value = 5.768654000000001
$ws.Cells.Item(11, 15).Value = 0.0888372291002696
$ws.Cells.Item(11, 16).Value = 0.0888372291002696
$ws.Cells.Item(11, 17).Value = 0.9136682637900001
$ws.Cells.Item(11, 18).Value = 8.22301437411
$ws.Cells.Item(11, 19).Value = 0.020548709503279
$ws.Cells.Item(11, 20).Value = 0.020548709503279

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.475155
$ws.Cells.Item(12, 8).Value = 1.425465
$ws.Cells.Item(12, 9).Value = 0.2313074114466796
$ws.Cells.Item(12, 10).Value = 0.2313074114466796
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.108471
$ws.Cells.Item(12, 14).Value = 0.325413
$ws.Cells.Item(12, 15).Value = 0.005011357802566427
$ws.Cells.Item(12, 16).Value = 0.005011357802566428
$ws.Cells.Item(12, 17).Value = 0.051540538005
$ws.Cells.Item(12, 18).Value = 0.463864842045
$ws.Cells.Item(12, 19).Value = 0.001159164201144761
$ws.Cells.Item(12, 20).Value = 0.001159164201144761

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.475155
$ws.Cells.Item(13, 8).Value = 1.425465
$ws.Cells.Item(13, 9).Value = 0.2313074114466796
$ws.Cells.Item(13, 10).Value = 0.2313074114466796
$ws.Cells.Item(13, 13).Value = 8.21807054472
$ws.Cells.Item(13, 14).Value = 73.96263490248
$ws.Cells.Item(13, 15).Value = 0.1848271971277737
$ws.Cells.Item(13, 16).Value = 0.1848271971277737
$ws.Cells.Item(13, 17).Value = 8.21807054472
$ws.Cells.Item(13, 18).Value = 73.96263490248
$ws.Cells.Item(13, 19).Value = 0.1848271971277737
$ws.Cells.Item(13, 20).Value = 0.1848271971277737
